# Auto-generated Excel COM-interop script
# Applies "Sitemaps updates and some corrections" edits:
#  - Populates column B ("Address") with updated relative platform/marketplace/storefront
#    documentation paths replacing the old URLs, for rows 2-158.
#  - Updates the active selection to B3 (and removes any scrolled viewport offset).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "?"
$ws.Range("B3").Value = "platform/user-guide/"
$ws.Range("B4").Value = "platform/user-guide/glossary/"
$ws.Range("B5").Value = "platform/user-guide/modules-installation/"
$ws.Range("B6").Value = "platform/user-guide/platform-overview/"
$ws.Range("B7").Value = "platform/user-guide/search/"
$ws.Range("B8").Value = "platform/user-guide/cart/settings/"
$ws.Range("B9").Value = "platform/user-guide/catalog/add-new-catalog/"
$ws.Range("B10").Value = "deleted page"
$ws.Range("B11").Value = "deleted page"
$ws.Range("B12").Value = "platform/user-guide/catalog/managing-categories/"
$ws.Range("B13").Value = "platform/user-guide/catalog/managing-product-variations/"
$ws.Range("B14").Value = "platform/user-guide/catalog/managing-products/"
$ws.Range("B15").Value = "platform/user-guide/catalog/managing-properties/"
$ws.Range("B16").Value = "platform/user-guide/catalog/managing-search-index/"
$ws.Range("B17").Value = "platform/user-guide/catalog/overview/"
$ws.Range("B18").Value = "platform/user-guide/catalog/product-indexing/"
$ws.Range("B19").Value = "platform/user-guide/catalog/setting-product-availability/"
$ws.Range("B20").Value = "platform/user-guide/catalog/settings/"
$ws.Range("B21").Value = "deleted page"
$ws.Range("B22").Value = "platform/user-guide/catalog/add-new-catalog/#view-catalogs"
$ws.Range("B23").Value = "platform/user-guide/catalog/managing-product-variations/#view-product-variations"
$ws.Range("B24").Value = "platform/user-guide/catalog/managing-properties/#view-properties"
$ws.Range("B25").Value = "platform/user-guide/catalog-personalization/example/"
$ws.Range("B26").Value = "platform/user-guide/catalog-personalization/overview/"
$ws.Range("B27").Value = "platform/user-guide/catalog-personalization/settings/"
$ws.Range("B28").Value = "deleted page"
$ws.Range("B29").Value = "platform/user-guide/catalog-personalization/user-groups/"
$ws.Range("B30").Value = "platform/user-guide/contacts/managing-contacts/"
$ws.Range("B31").Value = "platform/user-guide/contacts/overview/"
$ws.Range("B32").Value = "platform/user-guide/contacts/settings/"
$ws.Range("B33").Value = "platform/user-guide/content/managing-blogs/"
$ws.Range("B34").Value = "platform/user-guide/content/managing-linklists/"
$ws.Range("B35").Value = "platform/user-guide/content/managing-pages/"
$ws.Range("B36").Value = "platform/user-guide/content/managing-themes/"
$ws.Range("B37").Value = "platform/user-guide/content/overview/"
$ws.Range("B38").Value = "platform/user-guide/content/settings/"
$ws.Range("B39").Value = "platform/user-guide/contracts/creating-and-terminating-contracts/"
$ws.Range("B40").Value = "platform/user-guide/contracts/managing-contract-customers/"
$ws.Range("B41").Value = "platform/user-guide/contracts/managing-contract-prices/"
$ws.Range("B42").Value = "platform/user-guide/contracts/overview/"
$ws.Range("B43").Value = "platform/user-guide/generic-export/advanced-filtering/"
$ws.Range("B44").Value = "platform/user-guide/generic-export/assigning-permissions/"
$ws.Range("B45").Value = "platform/user-guide/generic-export/exporting-data/"
$ws.Range("B46").Value = "platform/user-guide/generic-export/overview/"
$ws.Range("B47").Value = "platform/user-guide/generic-export/settings/"
$ws.Range("B48").Value = "platform/user-guide/integrations/overview/"
$ws.Range("B49").Value = "platform/user-guide/integrations/avalara/address-validation/"
$ws.Range("B50").Value = "platform/user-guide/integrations/avalara/orders-synchronization/"
$ws.Range("B51").Value = "platform/user-guide/integrations/avalara/overview/"
$ws.Range("B52").Value = "platform/user-guide/integrations/avalara/settings/"
$ws.Range("B53").Value = "platform/user-guide/integrations/avalara/tax-type-configuration/"
$ws.Range("B54").Value = "platform/user-guide/integrations/avalara/taxes-calculation/"
$ws.Range("B55").Value = "platform/user-guide/integrations/google-analytics/integration/"
$ws.Range("B56").Value = "platform/user-guide/integrations/google-analytics/overview/"
$ws.Range("B57").Value = "platform/user-guide/integrations/google-analytics/settings/"
$ws.Range("B58").Value = "platform/user-guide/integrations/hotjar/getting-started/"
$ws.Range("B59").Value = "platform/user-guide/integrations/hotjar/overview/"
$ws.Range("B60").Value = "platform/user-guide/inventory/managing-fulfillment-centers/"
$ws.Range("B61").Value = "platform/user-guide/inventory/managing-inventory/"
$ws.Range("B62").Value = "platform/user-guide/inventory/overview/"
$ws.Range("B63").Value = "platform/user-guide/inventory/settings/"
$ws.Range("B64").Value = "platform/user-guide/marketing/advertising-spot/"
$ws.Range("B65").Value = "platform/user-guide/marketing/combining-active-promotions/"
$ws.Range("B66").Value = "platform/user-guide/marketing/dynamic-content-overview/"
$ws.Range("B67").Value = "platform/user-guide/marketing/managing-content-items/"
$ws.Range("B68").Value = "platform/user-guide/marketing/managing-content-placeholders/"
$ws.Range("B69").Value = "platform/user-guide/marketing/managing-promotions/"
$ws.Range("B70").Value = "platform/user-guide/marketing/managing-published-content/"
$ws.Range("B71").Value = "platform/user-guide/marketing/overview/"
$ws.Range("B72").Value = "platform/user-guide/marketing/promotion-rules/"
$ws.Range("B73").Value = "platform/user-guide/marketing/promotions-overview/"
$ws.Range("B74").Value = "platform/user-guide/marketing/publish-conditions/"
$ws.Range("B75").Value = "platform/user-guide/marketing/settings/"
$ws.Range("B76").Value = "marketplace/user-guide/"
$ws.Range("B77").Value = "marketplace/user-guide/Operator-portal/master-catalog-taxonomy-management/"
$ws.Range("B78").Value = "marketplace/user-guide/Operator-portal/overview/"
$ws.Range("B79").Value = "marketplace/user-guide/Operator-portal/products-management/"
$ws.Range("B80").Value = "marketplace/user-guide/Operator-portal/Commission-fees-setup/dynamic-commission-fees/"
$ws.Range("B81").Value = "marketplace/user-guide/Operator-portal/Commission-fees-setup/overview/"
$ws.Range("B82").Value = "marketplace/user-guide/Operator-portal/Commission-fees-setup/static-commission-fees/"
$ws.Range("B83").Value = "marketplace/user-guide/Operator-portal/Vendors-management/vendor-management/"
$ws.Range("B84").Value = "marketplace/user-guide/Operator-portal/Vendors-management/vendor-onboarding/"
$ws.Range("B85").Value = "marketplace/user-guide/Vendor-portal/offers-management-by-vendor/"
$ws.Range("B86").Value = "marketplace/user-guide/Vendor-portal/overview/"
$ws.Range("B87").Value = "marketplace/user-guide/Vendor-portal/products-management-by-vendor/"
$ws.Range("B88").Value = "platform/user-guide/notifications/notification-layouts/"
$ws.Range("B89").Value = "platform/user-guide/notifications/notification-list/"
$ws.Range("B90").Value = "platform/user-guide/notifications/notification-log/"
$ws.Range("B91").Value = "platform/user-guide/notifications/notification-templates/"
$ws.Range("B92").Value = "platform/user-guide/notifications/overview/"
$ws.Range("B93").Value = "platform/user-guide/notifications/settings/"
$ws.Range("B94").Value = "deleted page"
$ws.Range("B95").Value = "platform/user-guide/order-management/indexation/"
$ws.Range("B96").Value = "platform/user-guide/order-management/main-objects/"
$ws.Range("B97").Value = "platform/user-guide/order-management/managing-documents/"
$ws.Range("B98").Value = "platform/user-guide/order-management/managing-returns/"
$ws.Range("B99").Value = "platform/user-guide/order-management/notifications/"
$ws.Range("B100").Value = "platform/user-guide/order-management/overview/"
$ws.Range("B101").Value = "platform/user-guide/order-management/permissions/"
$ws.Range("B102").Value = "platform/user-guide/order-management/sending-order-information-to-avatax/"
$ws.Range("B103").Value = "platform/user-guide/order-management/settings/"
$ws.Range("B104").Value = "platform/user-guide/order-management/tracking-order-changes/"
$ws.Range("B105").Value = "platform/user-guide/pricing/adding-new-assignment/"
$ws.Range("B106").Value = "platform/user-guide/pricing/creating-new-price-list/#add-and-edit-product-prices"
$ws.Range("B107").Value = "platform/user-guide/pricing/creating-new-price-list/"
$ws.Range("B108").Value = "platform/user-guide/pricing/example/"
$ws.Range("B109").Value = "platform/user-guide/pricing/export-functionality/"
$ws.Range("B110").Value = "platform/user-guide/pricing/managing-pricing-module-settings/"
$ws.Range("B111").Value = "platform/user-guide/pricing/overview/"
$ws.Range("B112").Value = "platform/user-guide/pricing/troubleshooting-guide/"
$ws.Range("B113").Value = "platform/user-guide/pricing/viewing-price-list-in-catalog/"
$ws.Range("B114").Value = "platform/user-guide/security/api-key/"
$ws.Range("B115").Value = "platform/user-guide/security/login-on-behalf/"
$ws.Range("B116").Value = "platform/user-guide/security/managing-users/"
$ws.Range("B117").Value = "platform/user-guide/security/overview/"
$ws.Range("B118").Value = "platform/user-guide/security/roles-and-permissions/"
$ws.Range("B119").Value = "platform/user-guide/sitemaps/configuring-sitemaps/"
$ws.Range("B120").Value = "platform/user-guide/sitemaps/overview/"
$ws.Range("B121").Value = "platform/user-guide/sitemaps/settings/"
$ws.Range("B122").Value = "platform/user-guide/store/adding-new-store/"
$ws.Range("B123").Value = "platform/user-guide/store/configuring-store/"
$ws.Range("B124").Value = "platform/user-guide/store/overview/"
$ws.Range("B125").Value = "platform/user-guide/store/settings/"
$ws.Range("B126").Value = "storefront/user-guide/"
$ws.Range("B127").Value = "storefront/user-guide/account/addresses/"
$ws.Range("B128").Value = "storefront/user-guide/account/checkout-defaults/"
$ws.Range("B129").Value = "storefront/user-guide/account/company-info/"
$ws.Range("B130").Value = "storefront/user-guide/account/company-members/"
$ws.Range("B131").Value = "storefront/user-guide/account/dashboard/"
$ws.Range("B132").Value = "storefront/user-guide/account/lists/"
$ws.Range("B133").Value = "storefront/user-guide/account/orders/"
$ws.Range("B134").Value = "storefront/user-guide/account/overview/"
$ws.Range("B135").Value = "storefront/user-guide/account/profile/"
$ws.Range("B136").Value = "storefront/user-guide/account/quote-requests/"
$ws.Range("B137").Value = "storefront/user-guide/navigation/homepage-layout/"
$ws.Range("B138").Value = "storefront/user-guide/navigation/product-page-layout/"
$ws.Range("B139").Value = "storefront/user-guide/registration_and_signing_in/create-account/"
$ws.Range("B140").Value = "storefront/user-guide/registration_and_signing_in/password-management/"
$ws.Range("B141").Value = "storefront/user-guide/registration_and_signing_in/sign-in/"
$ws.Range("B142").Value = "deleted page"
$ws.Range("B143").Value = "storefront/user-guide/shopping/bulk-orders/"
$ws.Range("B144").Value = "storefront/user-guide/shopping/checkout-process/"
$ws.Range("B145").Value = "storefront/user-guide/shopping/compare-products/"
$ws.Range("B146").Value = "storefront/user-guide/shopping/lists/"
$ws.Range("B147").Value = "storefront/user-guide/shopping/searching-for-products/"
$ws.Range("B148").Value = "storefront/user-guide/shopping/submit-quotes/"
$ws.Range("B149").Value = "deleted page"
$ws.Range("B150").Value = "platform/user-guide/tasks/overview/"
$ws.Range("B151").Value = "platform/user-guide/tasks/roles-permissions/"
$ws.Range("B152").Value = "platform/user-guide/tasks/settings/"
$ws.Range("B153").Value = "platform/user-guide/tasks/using-application/"
$ws.Range("B154").Value = "platform/user-guide/thumbnails/generating-thumbnails/"
$ws.Range("B155").Value = "platform/user-guide/thumbnails/overview/"
$ws.Range("B156").Value = "platform/user-guide/thumbnails/settings/"
$ws.Range("B157").Value = "platform/user-guide/thumbnails/thumbnail-options/"
$ws.Range("B158").Value = "platform/user-guide/thumbnails/using-thumbnails/"

# Update the worksheet selection/view to match the saved state (cursor on B3, scrolled to top)
$ws.Range("A1").Select() | Out-Null
$ws.Range("B3").Select() | Out-Null
